$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14 - this shifts the existing rows 14-56 down to 15-57
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new weekly price entry
$ws.Range("A14").Value = 2
$ws.Range("B14").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C14").Value = "Coquimbo"
$ws.Range("D14").Value = 44608
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = 100112032
$ws.Range("G14").Value = "Zapallo italiano"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 500
$ws.Range("K14").Value = 7000
$ws.Range("L14").Value = 8000
$ws.Range("M14").Value = 7500
$ws.Range("N14").Value = "$/caja 60 unidades"
$ws.Range("O14").Value = "Provincia de Limarí"
$ws.Range("P14").Value = 125
$ws.Range("Q14").Value = 60
$ws.Range("R14").Value = "Hortaliza"
